# Revert 'cards' to commit 95cda46ab8 (Jun 25)
#
# The merchant_status sheet previously only carried the "business" columns
# (mst_code, mst_labe, updated_at). This restores the surrounding Airbyte
# sync metadata columns that wrap them:
#   _airbyte_ab_id | _airbyte_emitted_at | mst_code | mst_labe |
#   _airbyte_additional_properties | source_file_path | updated_at

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for the two leading metadata columns -------------
# Before: A=mst_code  B=mst_labe  C=updated_at
# After:              C=mst_code  D=mst_labe  E=updated_at
$ws.Range("A:B").Insert()

# --- Step 2: make room for the two metadata columns between mst_labe and
#             updated_at -------------------------------------------------
# Before: C=mst_code D=mst_labe E=updated_at
# After:  C=mst_code D=mst_labe          G=updated_at
$ws.Range("E:F").Insert()

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "_airbyte_ab_id"
$ws.Range("B1").Value = "_airbyte_emitted_at"
$ws.Range("E1").Value = "_airbyte_additional_properties"
$ws.Range("F1").Value = "source_file_path"

# Match the bold / centered / bordered header style already used by the
# surviving header cells (C1, D1, G1) by copying their format across.
$ws.Range("C1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows --------------------------------------------------------------
$abIds = @(
    "db8e07eb-2635-4fa8-8ee2-bc569f5d499f",
    "c5c5ca4f-a6e4-48c0-8778-293a709898e0",
    "478faea1-1fdd-4270-a9e6-2900c0efa8a0",
    "e3d82403-cf33-44bf-b12d-caa7ff5b4ac8"
)

$sourceFilePath = "s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/MERCHANT_STATUS/2024_08_06_1722929004063_0.parquet"

for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 2

    # _airbyte_ab_id
    $ws.Cells.Item($row, 1).Value = $abIds[$i]

    # _airbyte_emitted_at — same timestamp-number-format style as updated_at (col G)
    $ws.Range("G$row").Copy()
    $ws.Range("B$row").PasteSpecial(-4122)
    $excel.CutCopyMode = 0
    $ws.Cells.Item($row, 2).Value = 45510.3079196875

    # source_file_path
    $ws.Cells.Item($row, 6).Value = $sourceFilePath

    # updated_at refresh (slightly different cached timestamp than before)
    $ws.Cells.Item($row, 7).Value = 45511.29524080525
}

$wb.Save()
